$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.273.91"
$ws.Range("E2").Value = "  +0.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.689.59"
$ws.Range("E3").Value = "  +1.44%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.25"
$ws.Range("E5").Value = "  +0.80%  "

$ws.Range("E6").Value = "  +4.33%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2696"
$ws.Range("E8").Value = "  +2.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06441"
$ws.Range("E9").Value = "  +1.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.10"
$ws.Range("E10").Value = "  +3.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07457"
$ws.Range("E11").Value = "  +1.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.700.81"
$ws.Range("E12").Value = "  +1.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.557"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5863"
$ws.Range("E14").Value = "  +2.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008553"
$ws.Range("E15").Value = "  +1.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.66"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.322.09"
$ws.Range("E17").Value = "  +0.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.971"
$ws.Range("E18").Value = "  +0.83%  "

$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("E20").Value = "  +0.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.56"
$ws.Range("E21").Value = "  +1.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.238"
$ws.Range("E22").Value = "  +1.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.31"
$ws.Range("E24").Value = "  +1.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.683"
$ws.Range("E25").Value = "  +0.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1242"
$ws.Range("E26").Value = "  +6.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.88"
$ws.Range("E27").Value = "  +0.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06668"
$ws.Range("E28").Value = "  +15.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.350"
$ws.Range("E29").Value = "  +4.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.330"
$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.604"
$ws.Range("E31").Value = "  +3.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.558"
$ws.Range("E32").Value = "  +1.73%  "

$ws.Range("E33").Value = "  +1.35%  "

$ws.Range("E34").Value = "  +2.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6213"
$ws.Range("E35").Value = "  +3.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.388"
$ws.Range("E36").Value = "  +1.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.708"
$ws.Range("E37").Value = "  +2.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.284"
$ws.Range("E38").Value = "  +5.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01617"
$ws.Range("E39").Value = "  +1.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.103.37"
$ws.Range("E40").Value = "  +2.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8795"
$ws.Range("E41").Value = "  +2.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.016"
$ws.Range("E42").Value = "  +0.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.88"
$ws.Range("E43").Value = "  +1.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.837.40"
$ws.Range("E44").Value = "  +1.43%  "

$ws.Range("E45").Value = "  +2.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.93"
$ws.Range("E46").Value = "  +2.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.149"
$ws.Range("E47").Value = "  +0.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -0.32%  "

$ws.Range("E49").Value = "  +1.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4297"
$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.030"
$ws.Range("E51").Value = "  +2.64%  "
